# Add Leave Card entry (VL(1-0-0) / SL(1-0-0)) into the leave card table on Sheet1.
# This inserts one new row into the period table (Table1), shifting every row
# below it down by one, then fills in the new leave-usage data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item(1)

# --- Insert a new physical row above the current row 31 (period 45200) ---
$ws.Rows.Item(31).Insert()

# Grow the table definition so it covers the new row (A8:K136 -> A8:K137)
$tbl.Resize($ws.Range("A8:K137"))

# The freshly inserted row comes back with generic default formatting;
# clone the look of a normal data row (row 32) onto it.
$ws.Range("A32:K32").Copy()
$ws.Range("A31:K31").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# REMARKS (column K) on a row that carries a date needs the date-formatted
# style already used on row 28's remark cell.
$ws.Range("K28").Copy()
$ws.Range("K31").PasteSpecial(-4122)
$ws.Range("K30").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# --- Row 29 (period 45139): record 1.25 day VL/SL usage ---
$ws.Range("C29").Value = 1.25
$ws.Range("G29").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- Row 30 (period 45170): VL(1-0-0) particular, 1.25 earned, 1 day used ---
$ws.Range("B30").Value = "VL(1-0-0)"
$ws.Range("C30").Value = 1.25
$ws.Range("D30").Value = 1
$ws.Range("G30").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("K30").Value = 45170

# --- Row 31 (new row): SL(1-0-0) particular, 1 day used, remarks date ---
$ws.Range("B31").Value = "SL(1-0-0)"
$ws.Range("G31").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"
$ws.Range("H31").Value = 1
$ws.Range("K31").Value = 45195

# --- Row 137 (new trailing table row): restore the calculated-column formula ---
$ws.Range("G137").Formula = "=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"""",Table1[[#This Row],[EARNED]])"

# --- Refresh the remembered selection so it points at the newly-edited cell ---
$ws.Range("K31").Select()

$wb.Save()
